# Add Test Data for Russia / Finland / Hungary markets.
# New sheets are created by copying the "Netherlands" sheet (same template:
# layout, styles, merged cells, column widths) and are inserted right after
# "Denmark", which is currently the last sheet in the workbook.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$template = $wb.Worksheets.Item("Netherlands")
$lastSheet = $wb.Worksheets.Item("Denmark")

# --- Russia -----------------------------------------------------------
[void]$template.Copy($null, $lastSheet)
$russia = $wb.Worksheets.Item("Netherlands (2)")
$russia.Name = "Russia"
$russia.Range("B2").Value = "NGC-2929/T2910"
$russia.Range("B4").Value = "Russia Market"

# --- Finland ------------------------------------------------------------
[void]$template.Copy($null, $russia)
$finland = $wb.Worksheets.Item("Netherlands (2)")
$finland.Name = "Finland"
$finland.Range("B2").Value = "NGC-3130/T2943"
$finland.Range("B4").Value = "Finland Market"

# --- Hungary --------------------------------------------------------------
[void]$template.Copy($null, $finland)
$hungary = $wb.Worksheets.Item("Netherlands (2)")
$hungary.Name = "Hungary"
$hungary.Range("B2").Value = "NGC-3104/T2992"
$hungary.Range("B4").Value = "Hungary Market"

# Restore the selection on the template sheet (it lost focus/selection
# because of the copy operations) and mirror it on the new sheets.
[void]$template.Range("A7:A12").Select()
[void]$russia.Range("A7:A12").Select()
[void]$finland.Range("A7:A12").Select()

# Hungary ends up as the active sheet/tab, with the cursor left on H14.
[void]$hungary.Activate()
[void]$hungary.Range("H14").Select()

Write-Host "Added Russia, Finland and Hungary sheets"
